# updated main GSC export data
# Appends one more day (2025-12-30) to the "Chart" sheet's daily export table.
#
# The new date label must land as literal TEXT (a shared string), exactly like
# the 84 existing date rows above it - NOT as a real Excel date serial, which
# is what a plain `Range.Value = "2025-12-30"` assignment would produce (Excel
# auto-detects the "looks like a date" string and converts it). To keep the
# cell a plain string (same cell style as its neighbours), the new label is
# built via a formula (which is never subject to that auto-detection), copied,
# and pasted as a value into the target cell - then the scratch cell is
# removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$newRow = 86

# Stage the literal date text on a scratch cell far outside the used range.
# A formula result is always treated as text/number per its actual type, and
# is never reinterpreted as a "typed-in" date the way a direct Value/Formula
# string literal would be.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""2025-12-30"""

# Copy -> paste-as-value into the new row's date cell. Paste-values carries
# over only the literal content (as plain text), not the source formula and
# not any formatting, so the destination keeps its default (General) style -
# matching every other cell in the column.
$scratch.Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)  # xlPasteValues

# Clear the scratch cell so it leaves no trace in the saved sheet.
$scratch.Clear()

# Same row/day numbers as the previous row (no non-HTTPS URLs, 28 HTTPS URLs).
$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 28
